$d = $word.ActiveDocument

# New profit list entries to append after item 14 (Arancini)
$newLines = @(
    "15) ggddg profit: 1.50 ₪/min",
    "16) crabby patty profit: 0.03 ₪/min",
    "17) dddd profit: 0.02 ₪/min"
)

foreach ($line in $newLines) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $d.Paragraphs.Last.Range.Text = $line
}
